$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- String cell updates (order matters: it controls the order new shared
#     strings get appended to xl/sharedStrings.xml, which must match the
#     target file's table: ..., D3, CT21, C20, CT20) ---
$ws.Range("A2").Value = "CT21"
$ws.Range("B4").Value = "CT21"
$ws.Range("A8").Value = "CT21"
$ws.Range("B9").Value = "C20"
$ws.Range("B2").Value = "CT20"

# --- Clear cells that become blank ---
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D9").Value = ""

# --- Numeric cell updates (azimuth/distance correction) ---
$ws.Range("C3").Value = 2364556
$ws.Range("D3").Value = 70.811
$ws.Range("C5").Value = 2813648
$ws.Range("D5").Value = 65.651
$ws.Range("C7").Value = 3063251
$ws.Range("D7").Value = 86.345
$ws.Range("C9").Value = 750445

# --- Sheet view: drop scrolled top-left cell, move selection to D9 ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D9").Select()

# --- Workbook window size/position (best-effort; mirrors the author's
#     resized/moved Excel window from the saved workbookView) ---
$win = $wb.Windows.Item(1)
$win.Top = 2200
$win.Left = 2200
$win.Width = 14400
$win.Height = 7850
